# Auto-generated cell updates for Leviathan_Profits workbook
# (scheduled-runner style refresh of cached market-board values)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 244889.66
$ws.Range("J17").Value = 269238.53
$ws.Range("L17").Value = 807715.5900000001
$ws.Range("N17").Value = -808051.5900000001
$ws.Range("H134").Value = 120019
$ws.Range("J134").Value = 90279.5
$ws.Range("L134").Value = 90279.5
$ws.Range("N134").Value = -100419.5
$ws.Range("H137").Value = 1531
$ws.Range("I137").Value = 1355.3334
$ws.Range("J137").Value = 1636.4
$ws.Range("K137").Value = 4066.0002
$ws.Range("L137").Value = 4909.200000000001
$ws.Range("M137").Value = -1516.0002
$ws.Range("N137").Value = -10009.2
$ws.Range("H138").Value = 1767.4517
$ws.Range("J138").Value = 2229.6667
$ws.Range("L138").Value = 6689.000100000001
$ws.Range("N138").Value = -16969.0001
$ws.Range("H141").Value = 2768.8462
$ws.Range("I141").Value = 2926.818
$ws.Range("J141").Value = 1900
$ws.Range("K141").Value = 8780.454000000002
$ws.Range("L141").Value = 5700
$ws.Range("M141").Value = -3600.454000000002
$ws.Range("N141").Value = -16060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 907.0833
$ws.Range("I2").Value = 826.5
$ws.Range("J2").Value = 1793.5
$ws.Range("K2").Value = 826.5
$ws.Range("L2").Value = 1793.5
$ws.Range("M2").Value = -713.5
$ws.Range("N2").Value = -2019.5
$ws.Range("H32").Value = 7964.7036
$ws.Range("I32").Value = 7593.136
$ws.Range("K32").Value = 7593.136
$ws.Range("M32").Value = -7306.136
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 1384.963
$ws.Range("I74").Value = 1384.963
$ws.Range("K74").Value = 1384.963
$ws.Range("M74").Value = -510.963
$ws.Range("H77").Value = 1384.963
$ws.Range("I77").Value = 1384.963
$ws.Range("K77").Value = 6924.815
$ws.Range("M77").Value = -2556.815
$ws.Range("H97").Value = 1681.6
$ws.Range("I97").Value = 1681.6
$ws.Range("K97").Value = 1681.6
$ws.Range("M97").Value = -1185.6
$ws.Range("H116").Value = 907.0833
$ws.Range("I116").Value = 826.5
$ws.Range("J116").Value = 1793.5
$ws.Range("K116").Value = 826.5
$ws.Range("L116").Value = 1793.5
$ws.Range("M116").Value = 1467.5
$ws.Range("N116").Value = -6381.5
$ws.Range("H118").Value = 92124.75
$ws.Range("J118").Value = 92124.75
$ws.Range("L118").Value = 92124.75
$ws.Range("N118").Value = -95438.75
$ws.Range("H132").Value = 1648.2153
$ws.Range("I132").Value = 1603.4482
$ws.Range("K132").Value = 4810.3446
$ws.Range("M132").Value = -2280.3446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 89997
$ws.Range("J2").Value = 89997
$ws.Range("L2").Value = 89997
$ws.Range("N2").Value = -90223
$ws.Range("H3").Value = 907.0833
$ws.Range("I3").Value = 826.5
$ws.Range("J3").Value = 1793.5
$ws.Range("K3").Value = 826.5
$ws.Range("L3").Value = 1793.5
$ws.Range("M3").Value = -712.5
$ws.Range("N3").Value = -2021.5
$ws.Range("H82").Value = 19874.875
$ws.Range("I82").Value = 14749.75
$ws.Range("J82").Value = 25000
$ws.Range("K82").Value = 14749.75
$ws.Range("L82").Value = 25000
$ws.Range("M82").Value = -14366.75
$ws.Range("N82").Value = -25766
$ws.Range("H85").Value = 19874.875
$ws.Range("I85").Value = 14749.75
$ws.Range("J85").Value = 25000
$ws.Range("K85").Value = 14749.75
$ws.Range("L85").Value = 25000
$ws.Range("M85").Value = -13423.75
$ws.Range("N85").Value = -27652
$ws.Range("H94").Value = 1743.7693
$ws.Range("I94").Value = 1859.909
$ws.Range("J94").Value = 1105
$ws.Range("K94").Value = 1859.909
$ws.Range("L94").Value = 1105
$ws.Range("M94").Value = -1408.909
$ws.Range("N94").Value = -2007
$ws.Range("H105").Value = 14758.75
$ws.Range("I105").Value = 14758.75
$ws.Range("K105").Value = 14758.75
$ws.Range("M105").Value = -13011.75
$ws.Range("H134").Value = 954.6229
$ws.Range("I134").Value = 934.7719
$ws.Range("J134").Value = 1237.5
$ws.Range("K134").Value = 2804.3157
$ws.Range("L134").Value = 3712.5
$ws.Range("M134").Value = -269.3157000000001
$ws.Range("N134").Value = -8782.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12810.464
$ws.Range("I31").Value = 1860.5
$ws.Range("J31").Value = 78510.25
$ws.Range("K31").Value = 1860.5
$ws.Range("L31").Value = 78510.25
$ws.Range("M31").Value = -1565.5
$ws.Range("N31").Value = -79100.25
$ws.Range("H34").Value = 12810.464
$ws.Range("I34").Value = 1860.5
$ws.Range("J34").Value = 78510.25
$ws.Range("K34").Value = 1860.5
$ws.Range("L34").Value = 78510.25
$ws.Range("M34").Value = -1658.5
$ws.Range("N34").Value = -78914.25
$ws.Range("H58").Value = 2382
$ws.Range("I58").Value = 1456
$ws.Range("K58").Value = 1456
$ws.Range("M58").Value = -1253
$ws.Range("H120").Value = 22898
$ws.Range("J120").Value = 22898
$ws.Range("L120").Value = 22898
$ws.Range("N120").Value = -30156
$ws.Range("H132").Value = 3814.682
$ws.Range("I132").Value = 4664.5625
$ws.Range("J132").Value = 1548.3334
$ws.Range("K132").Value = 13993.6875
$ws.Range("L132").Value = 4645.0002
$ws.Range("M132").Value = -11463.6875
$ws.Range("N132").Value = -9705.0002
$ws.Range("H136").Value = 2382
$ws.Range("I136").Value = 1456
$ws.Range("K136").Value = 4368
$ws.Range("M136").Value = -1818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 132.2
$ws.Range("I40").Value = 132.2
$ws.Range("K40").Value = 528.8
$ws.Range("M40").Value = -459.8
$ws.Range("H70").Value = 5244
$ws.Range("I70").Value = 3995
$ws.Range("J70").Value = 5868.5
$ws.Range("K70").Value = 11985
$ws.Range("L70").Value = 17605.5
$ws.Range("M70").Value = -11670
$ws.Range("N70").Value = -18235.5
$ws.Range("H73").Value = 5244
$ws.Range("I73").Value = 3995
$ws.Range("J73").Value = 5868.5
$ws.Range("K73").Value = 11985
$ws.Range("L73").Value = 17605.5
$ws.Range("M73").Value = -10893
$ws.Range("N73").Value = -19789.5
$ws.Range("H105").Value = 9999.5
$ws.Range("I105").Value = 5000
$ws.Range("K105").Value = 15000
$ws.Range("M105").Value = -12379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 20699.46
$ws.Range("J57").Value = 21666.084
$ws.Range("L57").Value = 21666.084
$ws.Range("N57").Value = -23306.084
$ws.Range("H80").Value = 3622.0715
$ws.Range("I80").Value = 2673.4
$ws.Range("K80").Value = 2673.4
$ws.Range("M80").Value = -1675.4
$ws.Range("H83").Value = 3622.0715
$ws.Range("I83").Value = 2673.4
$ws.Range("K83").Value = 13367
$ws.Range("M83").Value = -8375
$ws.Range("H102").Value = 2642.52
$ws.Range("I102").Value = 2742.5217
$ws.Range("K102").Value = 2742.5217
$ws.Range("M102").Value = -1120.5217
$ws.Range("H123").Value = 61088.5
$ws.Range("J123").Value = 61088.5
$ws.Range("L123").Value = 61088.5
$ws.Range("N123").Value = -65988.5
$ws.Range("H126").Value = 1791.2
$ws.Range("I126").Value = 1464.375
$ws.Range("K126").Value = 4393.125
$ws.Range("M126").Value = -1923.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 47000
$ws.Range("J116").Value = 47000
$ws.Range("L116").Value = 47000
$ws.Range("N116").Value = -56178
$ws.Range("H133").Value = 119997.336
$ws.Range("J133").Value = 119997.336
$ws.Range("L133").Value = 119997.336
$ws.Range("N133").Value = -125057.336
$ws.Range("H136").Value = 4973.625
$ws.Range("I136").Value = 4358.8
$ws.Range("K136").Value = 13076.4
$ws.Range("M136").Value = -10526.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 27399
$ws.Range("I62").Value = 8497.5
$ws.Range("J62").Value = 40000
$ws.Range("K62").Value = 8497.5
$ws.Range("L62").Value = 40000
$ws.Range("M62").Value = -7873.5
$ws.Range("N62").Value = -41248
$ws.Range("H65").Value = 27399
$ws.Range("I65").Value = 8497.5
$ws.Range("J65").Value = 40000
$ws.Range("K65").Value = 42487.5
$ws.Range("L65").Value = 200000
$ws.Range("M65").Value = -39367.5
$ws.Range("N65").Value = -206240
$ws.Range("H117").Value = 48409
$ws.Range("J117").Value = 48409
$ws.Range("L117").Value = 48409
$ws.Range("N117").Value = -57587
$ws.Range("H120").Value = 51773
$ws.Range("J120").Value = 51773
$ws.Range("L120").Value = 51773
$ws.Range("N120").Value = -61449
